$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), styled like the existing headers (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New numeric data for columns I and J, rows 2-22
$values = @(
    @(5, 6),
    @(3, 5),
    @(7, 8),
    @(3, 6),
    @(7, 8),
    @(7, 8),
    @(6, 8),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(7, 9),
    @(1, 7),
    @(1, 4),
    @(1, 5),
    @(7, 8),
    @(1, 5),
    @(1, 4),
    @(5, 7),
    @(3, 5),
    @(1, 2),
    @(1, 1)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
